$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("s_curve")

# Rename header in A1 from "building_type" to "building_category"
$ws.Range("A1").Value = "building_category"

# Widen column A to fit the new, longer header text
# (ColumnWidth is expressed in characters; the stored XML column width adds
# the default 5/6-character padding Excel uses, so back that out here.)
$ws.Columns.Item(1).ColumnWidth = 16.1666666666667

# Update the saved window position/size metadata for the workbook
$win = $wb.Windows.Item(1)
$win.Left = -120
$win.Top = -120
$win.Width = 29040
$win.Height = 15840
